$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Unveiling Superconductivity's Potential" "Biology - Unveiling the Marvels of Life"

# Author
Replace-Text " Samantha Reynolds" " Alice Whitman"

# Email (user part only, keep ".edu")
Replace-Text "reynolds@quantumphysics" "awhitman@centralhigh"

# Body paragraph, sentence 1
Replace-Text "In the realm of physics, superconductivity stands as a captivating phenomenon, where certain materials exhibit the remarkable ability to conduct electricity without encountering resistance" "In the vast tapestry of natural wonders, biology occupies a unique place, unraveling the intricate mechanisms that govern the vibrant realm of life"

# Body paragraph, sentence 2
Replace-Text " This extraordinary characteristic unfolds at exceedingly low temperatures, often hovering near absolute zero" " This captivating subject unveils the enigmatic secrets of living organisms, tracing the symphony of life from birth to death"

# Body paragraph, sentences 3-4 merged into one
Replace-Text " Ever since its discovery in 1911, superconductivity has piqued the curiosity and ignited the imaginations of scientists and engineers alike. Its implications are profound, harboring the promise to revolutionize various fields, from power transmission and energy storage to computing and medical diagnostics" " Biology unveils the intricate ballet of cells, revealing the profound impact of molecules and genes that shape our existence"

# After first <br/><br/>
Replace-Text "With the advent of high-temperature superconductors in the 1980s, the practical applications of superconductivity moved closer to reality" "Through the lens of biology, we delve into the enigmatic enigma of heredity, tracing the intricate dance of chromosomes and DNA that define our traits"

Replace-Text " These materials, capable of exhibiting superconductivity at temperatures considerably higher than their conventional counterparts, opened up new avenues for exploration and development" " We explore the remarkable symphony of ecosystems, unraveling the intricate interconnectedness of life forms within their intricate ballet"

# Merge 4 sentences into one
Replace-Text " The potential benefits of superconductivity are immense. Imagine a world where electricity can be transmitted over long distances with minimal losses, paving the way for cleaner and more efficient power grids. Envision medical imaging devices with unprecedented sensitivity, enabling early detection and targeted treatment of diseases. Contemplate computers with lightning-fast processing speeds, empowering groundbreaking advances in artificial intelligence and machine learning" " From the delicate balance of predator and prey to the subtle interplay of symbiotic relationships, biology unveils the profound beauty and complexity of our natural world"

# After second <br/><br/>
Replace-Text "Superconductivity also holds the key to transformative technologies like magnetic levitation trains, levitating effortlessly above tracks and reaching breathtaking speeds" "Finally, biology empowers us to unravel the mysteries of human health, revealing the intricate mechanisms behind diseases and the remarkable resilience of the immune system"

# Merge 3 sentences into one
Replace-Text " These trains have the potential to revolutionize transportation, reducing travel times and energy consumption. Moreover, superconductivity could empower novel particle accelerators, unveiling the deepest secrets of the universe at even higher energies. Its applications extend far beyond these examples, touching diverse fields such as fusion energy, quantum computing, and materials science" " It holds the promise of unlocking the enigmatic potential of modern medicine, from groundbreaking vaccines to cutting-edge therapies, biology stands at the forefront of our quest to enhance human well-being"

# Summary heading paragraph - first sentence
Replace-Text "Unveiling the potential of superconductivity unveils a world of possibilities" "Biology unravels the intricate tapestry of life, unveiling the secrets of cellular processes, heredity, and ecosystems"

# Merge sentence + "research" (removes lastRenderedPageBreak run)
Replace-Text " This extraordinary property holds the promise to revolutionize energy, transportation, medicine, and scientific research" " It empowers us to comprehend the profound beauty of our natural world while simultaneously addressing the complexities of human health"

# Merge final two sentences
Replace-Text " As we continue to unlock the secrets of superconductivity, we stand on the threshold of a new era of innovation and technological advancements that have the power to reshape our world. Superconductivity's transformative potential is undeniable, and its journey towards practical applications is poised to redefine the boundaries of what is possible" " Through the lens of biology, we gain a deeper appreciation for the fragility and resilience of life, driving us toward a more compassionate and sustainable relationship with our planet and all of its inhabitants"

# Add a new empty paragraph at the end of the document (before sectPr)
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
